$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 11; existing rows 11-32 shift down to 12-33,
# carrying their formatting (incl. the date-style on column D) with them.
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11 with a new weekly price record
# (constant columns A,B,C,E,F,G,H,I,J,K,Q,R,T copied from the surrounding rows;
# D/L/M/N/O/P/S hold the new observation's data).
$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(11, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(11, 4).Value = 44690
$ws.Cells.Item(11, 5).Value = 15
$ws.Cells.Item(11, 6).Value = "Fruta"
$ws.Cells.Item(11, 7).Value = 100108
$ws.Cells.Item(11, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(11, 9).Value = 100108001
$ws.Cells.Item(11, 10).Value = "Guayaba"
$ws.Cells.Item(11, 11).Value = "Sin especificar"
$ws.Cells.Item(11, 12).Value = "Primera"
$ws.Cells.Item(11, 13).Value = 100
$ws.Cells.Item(11, 14).Value = 1600
$ws.Cells.Item(11, 15).Value = 1700
$ws.Cells.Item(11, 16).Value = 1650
$ws.Cells.Item(11, 17).Value = "$/kilo (en caja de 10 kilos )"
$ws.Cells.Item(11, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(11, 19).Value = 1650
$ws.Cells.Item(11, 20).Value = 1
